# edit.ps1
# Applies the "Kick out short reviews." commit to the document:
#  1. Collapse the split "Jeff" run (with proofErr gramStart/gramEnd wrapper)
#     in the arxiv BibTeX citation back into the surrounding run.
#  2. Mark the two inline result screenshots as NoProofing (w:noProof).
#  3. Replace the two trailing empty paragraphs at the end of the document
#     with the new Hebrew write-up (class-weighting note, cross-domain test
#     note, and the review/article success-rate bullet list).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge "Wu, " + "Jeff" + " and Child, ..." runs into a single run and
#    drop the now-unnecessary grammar-check proofErr markers around "Jeff".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Wu, Jeff and Child, Rewon",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Wu, Jeff and Child, Rewon", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Flag the two inline pictures (the train/test-on-articles screenshot and
#    the train-on-articles/test-on-reviews screenshot) as NoProofing.
# ---------------------------------------------------------------------------
$d.InlineShapes.Item(1).Range.NoProofing = $true
$d.InlineShapes.Item(2).Range.NoProofing = $true

# ---------------------------------------------------------------------------
# 3) Replace the final two (empty) paragraphs with the new write-up.
#    The next-to-last paragraph is replaced by thirteen new paragraphs;
#    the very last paragraph (the document must always end on a paragraph
#    mark) is replaced in place by the fourteenth/final new paragraph so no
#    stray empty paragraph is left behind before the sectPr.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($count - 1)

$block1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
      </w:pPr>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="5"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">הוספנו </w:t>
      </w:r>
      <w:r>
        <w:t>class weighting</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> ל-</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>SVM</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="5"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>כדי לבדוק עוד: הרצנו את הבדיקה לא ב-</w:t>
      </w:r>
      <w:r>
        <w:t>cross-domain</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">: פעם אחת אימון על מאמרים ובדיקה על מאמרים, ופעם אחת אימון על </w:t>
      </w:r>
      <w:r>
        <w:t>reviews</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> ובדיקה על </w:t>
      </w:r>
      <w:r>
        <w:t>reviews</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">אחוזי </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>הצלחה בשלב הזה:</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>review-review</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">בלי </w:t>
      </w:r>
      <w:r>
        <w:t>perplexity</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> 70%</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">עם </w:t>
      </w:r>
      <w:r>
        <w:t>perplexity</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> 58%</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">עם </w:t>
      </w:r>
      <w:r>
        <w:t>perplexity</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">, עם 7 מחלקות במקום 9 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> 65%</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>מאמרים-מאמרים:</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">בלי </w:t>
      </w:r>
      <w:r>
        <w:t>perplexity</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> 74% </w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">עם </w:t>
      </w:r>
      <w:r>
        <w:t>perplexity</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>66%</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>מאמרים-</w:t>
      </w:r>
      <w:r>
        <w:t>reviews</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">בלי </w:t>
      </w:r>
      <w:r>
        <w:t>perplexity</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> 24%</w:t>
      </w:r>
    </w:p>
'@
$secondToLast.Range.InsertXML($block1)

$newCount = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($newCount)

$block2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">עם </w:t>
      </w:r>
      <w:r>
        <w:t>perplexity</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> 30% אבל על כולם חוץ מאחד הוא חזה </w:t>
      </w:r>
      <w:r>
        <w:t>Richard riley</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
'@
$last.Range.InsertXML($block2)
